# The stock report has several items that appear as two adjacent rows
# (same item description) whose Item Code / Rate / Qty / Value figures had
# been swapped between the two rows. This corrects each such pair by
# swapping columns B (Item Code), E (Rate), F (Qty) and G (Value) between
# the two adjacent rows so the figures line up with the correct row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    112,
    127,
    192,
    229,
    232,
    243,
    364,
    372,
    380,
    382,
    385,
    473,
    572
)

foreach ($r in $rowPairs) {
    $r2 = $r + 1

    $bTop = $ws.Cells.Item($r, 2).Value2
    $eTop = $ws.Cells.Item($r, 5).Value2
    $fTop = $ws.Cells.Item($r, 6).Value2
    $gTop = $ws.Cells.Item($r, 7).Value2

    $bBot = $ws.Cells.Item($r2, 2).Value2
    $eBot = $ws.Cells.Item($r2, 5).Value2
    $fBot = $ws.Cells.Item($r2, 6).Value2
    $gBot = $ws.Cells.Item($r2, 7).Value2

    $ws.Cells.Item($r, 2).Value = $bBot
    $ws.Cells.Item($r, 5).Value = $eBot
    $ws.Cells.Item($r, 6).Value = $fBot
    $ws.Cells.Item($r, 7).Value = $gBot

    $ws.Cells.Item($r2, 2).Value = $bTop
    $ws.Cells.Item($r2, 5).Value = $eTop
    $ws.Cells.Item($r2, 6).Value = $fTop
    $ws.Cells.Item($r2, 7).Value = $gTop
}
